$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ol-Liq")
$ws.Activate()

# Replace the "Reference" text for rows 8-10 with new imaginary sample labels
$ws.Range("B8").Value = "Imaginary water-rich sample 1"
$ws.Range("B9").Value = "Imaginary water-rich sample 2"
$ws.Range("B10").Value = "Imaginary water-rich sample 3"

# Update the H2O_Liq values for the same rows
$ws.Range("N8").Value = 2
$ws.Range("N9").Value = 4
$ws.Range("N10").Value = 6

# Move the selection to N11, mirroring where the user last clicked
$ws.Range("N11").Select()
